$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = '[''Romania'', 2, -1, 2]'
$ws.Cells.Item(2, 15).Value = '[]'
$ws.Cells.Item(4, 8).Value = '[''Slovakia'', 4, 0, 3]'
$ws.Cells.Item(4, 13).Value = '[''Slovakia'', ''Northern Ireland'', ''Albania'', ''Portugal'']'
$ws.Cells.Item(4, 16).Value = 1
$ws.Cells.Item(4, 17).Value = 2
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(8, 9).Value = '[''Northern Ireland'', 4, 1, 2]'
$ws.Cells.Item(8, 13).Value = '[''Northern Ireland'', ''Slovakia'', ''Albania'', ''Portugal'']'
$ws.Cells.Item(8, 16).Value = 1
$ws.Cells.Item(8, 17).Value = 3
$ws.Cells.Item(9, 16).Value = 1
$ws.Cells.Item(9, 17).Value = 4
$ws.Cells.Item(10, 17).Value = 4
$ws.Cells.Item(11, 10).Value = '[''Czech Republic'', 2, -1, 2]'
$ws.Cells.Item(11, 17).Value = 4
$ws.Cells.Item(12, 17).Value = 4
$ws.Cells.Item(13, 17).Value = 5
$ws.Cells.Item(14, 17).Value = 5
$ws.Cells.Item(15, 17).Value = 6
$ws.Cells.Item(16, 17).Value = 6
$ws.Cells.Item(17, 12).Value = '[''Portugal'', 3, 0, 1]'
$ws.Cells.Item(17, 13).Value = '[''Slovakia'', ''Northern Ireland'', ''Portugal'', ''Turkey'']'
$ws.Cells.Item(17, 14).Value = '[''Albania'', ''Sweden'']'
$ws.Cells.Item(17, 16).Value = 1
$ws.Cells.Item(17, 17).Value = 7
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 7
$ws.Cells.Item(19, 17).Value = 8
$ws.Cells.Item(20, 17).Value = 9
$ws.Cells.Item(21, 17).Value = 10
$ws.Cells.Item(22, 17).Value = 11
$ws.Cells.Item(23, 17).Value = 12
$ws.Cells.Item(24, 17).Value = 12
$ws.Cells.Item(25, 17).Value = 13
$ws.Cells.Item(26, 17).Value = 14
$ws.Cells.Item(27, 11).Value = '[''Sweden'', 2, -1, 1]'
$ws.Cells.Item(27, 17).Value = 14
$ws.Cells.Item(28, 17).Value = 14
$ws.Cells.Item(29, 17).Value = 15
$ws.Cells.Item(30, 7).Value = '[''Switzerland'', 2, -3, 1]'
$ws.Cells.Item(30, 14).Value = '[''Switzerland'', ''Croatia'']'
$ws.Cells.Item(36, 9).Value = '[''Austria'', 4, 0, 3]'
$ws.Cells.Item(36, 13).Value = '[''Austria'', ''Switzerland'', ''Portugal'', ''Finland'']'
$ws.Cells.Item(36, 16).Value = 1
$ws.Cells.Item(36, 17).Value = 2
$ws.Cells.Item(37, 17).Value = 3
$ws.Cells.Item(38, 17).Value = 3
$ws.Cells.Item(39, 17).Value = 3
$ws.Cells.Item(40, 17).Value = 3
$ws.Cells.Item(41, 8).Value = '[''Finland'', 4, 0, 1]'
$ws.Cells.Item(41, 13).Value = '[''Finland'', ''Switzerland'', ''Portugal'', ''Ukraine'']'
$ws.Cells.Item(41, 16).Value = 1
$ws.Cells.Item(41, 17).Value = 4
$ws.Cells.Item(42, 17).Value = 5
$ws.Cells.Item(43, 17).Value = 6
$ws.Cells.Item(44, 17).Value = 7
$ws.Cells.Item(45, 17).Value = 8
$ws.Cells.Item(46, 17).Value = 9
$ws.Cells.Item(47, 17).Value = 9
$ws.Cells.Item(48, 17).Value = 9
$ws.Cells.Item(49, 10).Value = '[''Croatia'', 2, -1, 1]'
$ws.Cells.Item(49, 17).Value = 9
$ws.Cells.Item(50, 17).Value = 9
$ws.Cells.Item(51, 17).Value = 10
$ws.Cells.Item(52, 17).Value = 11
$ws.Cells.Item(53, 17).Value = 12
$ws.Cells.Item(54, 17).Value = 13
$ws.Cells.Item(55, 11).Value = '[''Spain'', 3, 0, 1]'
$ws.Cells.Item(55, 13).Value = '[''Czech Republic'', ''Switzerland'', ''Portugal'', ''Spain'']'
$ws.Cells.Item(55, 14).Value = '[''Ukraine'', ''Finland'']'
$ws.Cells.Item(55, 16).Value = 1
$ws.Cells.Item(55, 17).Value = 14
$ws.Cells.Item(56, 16).Value = 0
$ws.Cells.Item(56, 17).Value = 14
$ws.Cells.Item(57, 17).Value = 15
$ws.Cells.Item(58, 17).Value = 15
$ws.Cells.Item(59, 17).Value = 15
$ws.Cells.Item(60, 17).Value = 15
$ws.Cells.Item(61, 17).Value = 15
$ws.Cells.Item(62, 17).Value = 15
$ws.Cells.Item(63, 17).Value = 15
$ws.Cells.Item(64, 17).Value = 15
$ws.Cells.Item(65, 17).Value = 15
$ws.Cells.Item(66, 12).Value = '[''Portugal'', 4, 1, 5]'
$ws.Cells.Item(66, 13).Value = '[''Portugal'', ''Czech Republic'', ''Switzerland'', ''Ukraine'']'
$ws.Cells.Item(66, 16).Value = 1
$ws.Cells.Item(66, 17).Value = 16
$ws.Cells.Item(67, 17).Value = 17
$ws.Cells.Item(68, 17).Value = 17
$ws.Cells.Item(69, 17).Value = 17
$ws.Cells.Item(70, 17).Value = 18
$ws.Cells.Item(71, 17).Value = 19
$ws.Cells.Item(72, 17).Value = 20
$ws.Cells.Item(73, 17).Value = 21
$ws.Cells.Item(74, 17).Value = 22
$ws.Cells.Item(75, 7).Value = '[''Scotland'', 2, -4, 2]'
$ws.Cells.Item(75, 13).Value = '[''Austria'', ''Belgium'', ''Denmark'', ''Scotland'']'
$ws.Cells.Item(75, 14).Value = '[''Albania'', ''Czech Republic'']'
$ws.Cells.Item(76, 16).Value = 0
$ws.Cells.Item(76, 17).Value = 0
$ws.Cells.Item(77, 17).Value = 0
$ws.Cells.Item(78, 17).Value = 1
$ws.Cells.Item(79, 8).Value = '[''Albania'', 2, -1, 3]'
$ws.Cells.Item(79, 17).Value = 1
$ws.Cells.Item(80, 17).Value = 1
$ws.Cells.Item(81, 17).Value = 2
$ws.Cells.Item(82, 17).Value = 3
$ws.Cells.Item(83, 10).Value = '[''Austria'', 4, 1, 3]'
$ws.Cells.Item(83, 17).Value = 3
$ws.Cells.Item(84, 17).Value = 4
$ws.Cells.Item(85, 17).Value = 5
$ws.Cells.Item(86, 17).Value = 5
$ws.Cells.Item(87, 17).Value = 6
$ws.Cells.Item(88, 17).Value = 7
$ws.Cells.Item(89, 17).Value = 7
$ws.Cells.Item(90, 17).Value = 8
$ws.Cells.Item(91, 11).Value = '[''Belgium'', 4, 1, 2]'
$ws.Cells.Item(91, 13).Value = '[''Belgium'', ''Netherlands'', ''Hungary'', ''Denmark'']'
$ws.Cells.Item(91, 16).Value = 1
$ws.Cells.Item(94, 12).Value = '[''Czech Republic'', 2, -1, 2]'
$ws.Cells.Item(94, 14).Value = '[''Czech Republic'', ''Croatia'']'
$ws.Cells.Item(97, 15).Value = '[''Georgia'', ''Netherlands'']'
$ws.Cells.Item(98, 15).Value = '[''Georgia'', ''Netherlands'']'
$ws.Cells.Item(99, 15).Value = '[''Georgia'', ''Netherlands'']'
$ws.Cells.Item(100, 9).Value = '[''Denmark'', 3, 0, 2]'
$ws.Cells.Item(100, 13).Value = '[''Netherlands'', ''Georgia'', ''Slovakia'', ''Denmark'']'
$ws.Cells.Item(100, 14).Value = '[''Hungary'', ''Croatia'']'
$ws.Cells.Item(100, 15).Value = '[''Georgia'', ''Netherlands'']'
$ws.Cells.Item(100, 16).Value = 1
$ws.Cells.Item(100, 17).Value = 14
